$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.912.51"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.876.25"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'0.7409"
$ws.Range("E5").Value = "  -4.02%  "
$ws.Range("D6").Value = "'242.63"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.9994"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'0.3164"
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("D9").Value = "'0.07214"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").Value = "'24.75"
$ws.Range("E10").Value = "  -3.05%  "
$ws.Range("D11").Value = "'0.08386"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").Value = "'0.7520"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.432"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.836.92"
$ws.Range("E14").Value = "  -7.72%  "
$ws.Range("D15").Value = "'92.80"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "29.891.24"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "'6.085"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D19").Value = "'13.61"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").Value = "'0.000007861"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").Value = "'0.9980"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "2.127.27"
$ws.Range("E22").Value = "  -5.41%  "
$ws.Range("D23").Value = "'8.060"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'0.1554"
$ws.Range("E25").Value = "  -5.05%  "
$ws.Range("D26").Value = "'9.280"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").Value = "'165.51"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").Value = "'18.68"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "'2.042"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").Value = "'1.501"
$ws.Range("E30").Value = "  +4.76%  "
$ws.Range("D31").Value = "'4.596"
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("D32").Value = "'1.537"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "'4.290"
$ws.Range("E33").Value = "  +4.44%  "
$ws.Range("D34").Value = "'0.05318"
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("D35").Value = "'1.240"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'0.7562"
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("D37").Value = "'1.003"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").Value = "'2.755"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").Value = "'0.4505"
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("D42").Value = "1.111.09"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "'6.061"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "'72.31"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("D45").Value = "'0.8566"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "'103.26"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").Value = "'7.642"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").Value = "'1.854"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.484"
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'2.944"
$ws.Range("E51").Value = "  -1.61%  "
